# -----------------------------------------------------------------------
# 1) Table on slide 16 (graphicFrame "Google Shape;213;p29") switches from
#    the deck's custom table style ("Table_0") to the built-in
#    "Medium Style 2 - Accent 1" table style.
# -----------------------------------------------------------------------
$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{4CFB2FE5-2CF5-4CB3-B049-E5CAC82050E7}")

# -----------------------------------------------------------------------
# 2) The presentation's theme reverts from the "Integral" design back to
#    the default "Office Theme" palette (dk1/lt1/dk2/lt2/accent1-6/
#    hlink/folHlink), applied through the slide master's theme color
#    scheme (MsoThemeColorSchemeIndex order 1-12).
# -----------------------------------------------------------------------
function ToBGRInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - the stock "Office" theme.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ToBGRInt($officeThemeColors[$i - 1])
}
